try {
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 49: fill in EARNED (VL) value ---
$ws.Range("C49").Value = 1.25

# --- Row 50: add SL(1-0-0) entry, earned value, days=1, end date ---
$ws.Range("B50").Value = "SL(1-0-0)"
$ws.Range("C50").Value = 1.25
$ws.Range("H50").Value = 1
$ws.Range("K49").Copy()
$ws.Range("K50").PasteSpecial(-4122)
$ws.Range("K50").Value = 45208

# --- Row 51: add SL(1-0-0) entry, earned value, days=1, end date ---
$ws.Range("B51").Value = "SL(1-0-0)"
$ws.Range("C51").Value = 1.25
$ws.Range("H51").Value = 1
$ws.Range("K49").Copy()
$ws.Range("K51").PasteSpecial(-4122)
$ws.Range("K51").Value = 45239

# --- Insert a new row before (old) row 53 to add a "2024" year-section header ---
$ws.Rows("53:53").Insert()

# Restore formatting for the newly inserted blank row 53 by copying it
# from the row directly below (which now holds the old row 53 content/format).
$ws.Range("A54:K54").Copy()
$ws.Range("A53:K53").PasteSpecial(-4122)

# Apply the "year header" formatting (as used for the 2021/2022 headers) to A53.
$ws.Range("A21").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$ws.Range("A53").Value = "'2024"

# Grow Table1 to include the newly added row (was A8:K136, now A8:K137).
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A8:K137"))

# Fix the calculated column formula text/value for the brand new last row
# so it matches the long-form table formula used throughout the sheet.
$ws.Range("G137").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),`"`",Table1[[#This Row],[EARNED]])"

# Update the active cell selections to reflect where the user ended up editing.
$win = $excel.ActiveWindow
$win.Panes.Item(1).Activate() | Out-Null
$ws.Range("I9").Select() | Out-Null
$win.Panes.Item(2).Activate() | Out-Null
$ws.Range("F48").Select() | Out-Null

$wb.Application.CalculateFull()

Write-Output "OK"
} catch {
Write-Output "ERROR: $_"
}
